$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.02802633333333333
$ws.Range("H2").Value = 0.084079
$ws.Range("I2").Value = 0.09156020089470454
$ws.Range("J2").Value = 0.09156020089470451
$ws.Range("M2").Value = 13.71977066666667
$ws.Range("N2").Value = 41.159312
$ws.Range("O2").Value = 0.5515038136402627
$ws.Range("P2").Value = 0.5515038136402626
$ws.Range("Q2").Value = 0.3845148659608889
$ws.Range("R2").Value = 3.460633793648
$ws.Range("S2").Value = 0.05049579997109815
$ws.Range("T2").Value = 0.05049579997109812

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.02802633333333333
$ws.Range("H3").Value = 0.084079
$ws.Range("I3").Value = 0.09156020089470454
$ws.Range("J3").Value = 0.09156020089470451
$ws.Range("O3").Value = 0.172077867958883
$ws.Range("P3").Value = 0.1720778679588829
$ws.Range("Q3").Value = 0.1199746886541111
$ws.Range("R3").Value = 1.079772197887
$ws.Range("S3").Value = 0.01575548415984777
$ws.Range("T3").Value = 0.01575548415984776

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.02802633333333333
$ws.Range("H4").Value = 0.084079
$ws.Range("I4").Value = 0.09156020089470454
$ws.Range("J4").Value = 0.09156020089470451
$ws.Range("O4").Value = 0.2764183184008545
$ws.Range("P4").Value = 0.2764183184008545
$ws.Range("Q4").Value = 0.1927220628765555
$ws.Range("R4").Value = 1.734498565889
$ws.Range("S4").Value = 0.02530891676375864
$ws.Range("T4").Value = 0.02530891676375863

$ws.Range("I5").Value = 0.8397895222870286
$ws.Range("J5").Value = 0.8397895222870285
$ws.Range("M5").Value = 13.71977066666667
$ws.Range("N5").Value = 41.159312
$ws.Range("O5").Value = 0.5515038136402627
$ws.Range("P5").Value = 0.5515038136402626
$ws.Range("Q5").Value = 3.526767661518222
$ws.Range("R5").Value = 31.740908953664
$ws.Range("S5").Value = 0.4631471241964307
$ws.Range("T5").Value = 0.4631471241964305

$ws.Range("I6").Value = 0.8397895222870286
$ws.Range("J6").Value = 0.8397895222870285
$ws.Range("O6").Value = 0.172077867958883
$ws.Range("P6").Value = 0.1720778679588829
$ws.Range("S6").Value = 0.1445091905293607
$ws.Range("T6").Value = 0.1445091905293607

$ws.Range("I7").Value = 0.8397895222870286
$ws.Range("J7").Value = 0.8397895222870285
$ws.Range("O7").Value = 0.2764183184008545
$ws.Range("P7").Value = 0.2764183184008545
$ws.Range("S7").Value = 0.2321332075612373
$ws.Range("T7").Value = 0.2321332075612373

$ws.Range("I8").Value = 0.06865027681826696
$ws.Range("J8").Value = 0.06865027681826695
$ws.Range("M8").Value = 13.71977066666667
$ws.Range("N8").Value = 41.159312
$ws.Range("O8").Value = 0.5515038136402627
$ws.Range("P8").Value = 0.5515038136402626
$ws.Range("Q8").Value = 0.2883026875324445
$ws.Range("R8").Value = 2.594724187792
$ws.Range("S8").Value = 0.03786088947273395
$ws.Range("T8").Value = 0.03786088947273394

$ws.Range("I9").Value = 0.06865027681826696
$ws.Range("J9").Value = 0.06865027681826695
$ws.Range("O9").Value = 0.172077867958883
$ws.Range("P9").Value = 0.1720778679588829
$ws.Range("S9").Value = 0.01181319326967451
$ws.Range("T9").Value = 0.0118131932696745

$ws.Range("I10").Value = 0.06865027681826696
$ws.Range("J10").Value = 0.06865027681826695
$ws.Range("O10").Value = 0.2764183184008545
$ws.Range("P10").Value = 0.2764183184008545
$ws.Range("S10").Value = 0.01897619407585852
$ws.Range("T10").Value = 0.01897619407585851

